$wb = $excel.ActiveWorkbook

# --- Sheet "Matriz_Resultados": set corrected win/loss/tie cells to 0 ---
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$matrizZeroCells = @("C2","D2","E2","G2","H2","I2","B3","F3","G3","B4","F4","G4","H4","B5","F5","G5","J5","C6","D6","E6","G6","H6","I6","B7","C7","D7","E7","F7","J7","B8","D8","F8","B9","F9","J9","E10","G10","I10")
foreach ($addr in $matrizZeroCells) {
    $wsMatriz.Range($addr).Value = 0
}

# --- Sheet "P_valores": update corrected p-values ---
$wsP = $wb.Worksheets.Item("P_valores")
$pValueUpdates = @{
    "C2" = 0.00232213727036501
    "D2" = 0.003364046601608939
    "E2" = 0.01051944376908454
    "F2" = 0.04955890974999688
    "G2" = 0.003435303304379911
    "H2" = 0.00343121755398057
    "I2" = 0.002508875049209935
    "J2" = 0.02642537288565316
    "B3" = 0.00232213727036501
    "D3" = 0.000029170709975323561025106756
    "E3" = 0.0004163570426174434
    "F3" = 0.001819778583110443
    "G3" = 0.001461336250830003
    "H3" = 0.001053524306092068
    "I3" = 0.03735020695915825
    "J3" = 0.0006577743106987999
    "B4" = 0.003364046601608939
    "C4" = 0.000029170709975323561025106756
    "E4" = 0.0009375623716423309
    "F4" = 0.002700922160544961
    "G4" = 0.01616760701843534
    "H4" = 0.007008314228470036
    "I4" = 0.1985018123831537
    "J4" = 0.001092096716558189
    "B5" = 0.01051944376908454
    "C5" = 0.0004163570426174434
    "D5" = 0.0009375623716423309
    "F5" = 0.009319117563492574
    "G5" = 0.02078658305101411
    "H5" = 0.01204500379729545
    "I5" = 0.08614962132361881
    "J5" = 0.0081070354815167
    "B6" = 0.04955890974999688
    "C6" = 0.001819778583110443
    "D6" = 0.002700922160544961
    "E6" = 0.009319117563492574
    "G6" = 0.002639943580309945
    "H6" = 0.002652109936424951
    "I6" = 0.002060209066193197
    "J6" = 0.03490487857461577
    "B7" = 0.003435303304379911
    "C7" = 0.001461336250830003
    "D7" = 0.01616760701843534
    "E7" = 0.02078658305101411
    "F7" = 0.002639943580309945
    "H7" = 0.4327896701971903
    "I7" = 0.6510782794572643
    "J7" = 0.002084086427770915
    "B8" = 0.00343121755398057
    "C8" = 0.001053524306092068
    "D8" = 0.007008314228470036
    "E8" = 0.01204500379729545
    "F8" = 0.002652109936424951
    "G8" = 0.4327896701971903
    "I8" = 0.4984043110240182
    "J8" = 0.001358454403508125
    "B9" = 0.002508875049209935
    "C9" = 0.03735020695915825
    "D9" = 0.1985018123831537
    "E9" = 0.08614962132361881
    "F9" = 0.002060209066193197
    "G9" = 0.6510782794572643
    "H9" = 0.4984043110240182
    "J9" = 0.003348033413496765
    "B10" = 0.02642537288565316
    "C10" = 0.0006577743106987999
    "D10" = 0.001092096716558189
    "E10" = 0.0081070354815167
    "F10" = 0.03490487857461577
    "G10" = 0.002084086427770915
    "H10" = 0.001358454403508125
    "I10" = 0.003348033413496765
}
foreach ($addr in $pValueUpdates.Keys) {
    $wsP.Range($addr).Value = $pValueUpdates[$addr]
}

# --- Sheet "Estadisticos_DM": update corrected DM statistics ---
$wsDM = $wb.Worksheets.Item("Estadisticos_DM")
$dmStatUpdates = @{
    "C2" = 3.442660304579882
    "D2" = 3.287018271814858
    "E2" = 2.796473722301279
    "F2" = 2.078248836747459
    "G2" = 3.2781720493548
    "H2" = 3.278674436326821
    "I2" = 3.410297521038933
    "J2" = 2.379654476726118
    "B3" = -3.442660304579882
    "D3" = -5.244970318175312
    "E3" = -4.151892140990231
    "F3" = -3.544308304162992
    "G3" = -3.635376023382044
    "H3" = -3.770632368888417
    "I3" = -2.215914064081026
    "J3" = -3.964416740307003
    "B4" = -3.287018271814858
    "C4" = 5.244970318175312
    "E4" = -3.818694910841425
    "F4" = -3.379381633110052
    "G4" = -2.604939528684409
    "H4" = -2.973714205662922
    "I4" = -1.325820200266363
    "J4" = -3.755798272031607
    "B5" = -2.796473722301279
    "C5" = 4.151892140990231
    "D5" = 3.818694910841425
    "F5" = -2.849693640197735
    "G5" = 2.490659316536247
    "H5" = 2.736604789728622
    "I5" = 1.796567395782656
    "J5" = -2.91052663255958
    "B6" = -2.078248836747459
    "C6" = 3.544308304162992
    "D6" = 3.379381633110052
    "E6" = 2.849693640197735
    "G6" = 3.388959000717004
    "H6" = 3.38703101784544
    "I6" = 3.492627954080702
    "J6" = 2.248335654680929
    "B7" = -3.2781720493548
    "C7" = 3.635376023382044
    "D7" = 2.604939528684409
    "E7" = -2.490659316536247
    "F7" = -3.388959000717004
    "H7" = -0.7990752195082508
    "I7" = 0.4585172155200827
    "J7" = -3.487822442158532
    "B8" = -3.278674436326821
    "C8" = 3.770632368888417
    "D8" = 2.973714205662922
    "E8" = -2.736604789728622
    "F8" = -3.38703101784544
    "G8" = 0.7990752195082508
    "I8" = 0.6883871919935181
    "J8" = -3.665608963564621
    "B9" = -3.410297521038933
    "C9" = 2.215914064081026
    "D9" = 1.325820200266363
    "E9" = -1.796567395782656
    "F9" = -3.492627954080702
    "G9" = -0.4585172155200827
    "H9" = -0.6883871919935181
    "J9" = -3.289031301444221
    "B10" = -2.379654476726118
    "C10" = 3.964416740307003
    "D10" = 3.755798272031607
    "E10" = 2.91052663255958
    "F10" = -2.248335654680929
    "G10" = 3.487822442158532
    "H10" = 3.665608963564621
    "I10" = 3.289031301444221
}
foreach ($addr in $dmStatUpdates.Keys) {
    $wsDM.Range($addr).Value = $dmStatUpdates[$addr]
}

# --- Sheet "Resumen": rewrite summary rows (re-ranked after correction) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$resumenRows = @{
    2 = @("Sieve Bootstrap", 4, 0, 4, 50, 0.5800196476824288)
    3 = @("LSPM", 2, 1, 5, 25, 0.8081876453848597)
    4 = @("AV-MCPS", 1, 1, 6, 12.5, 1.243186337533633)
    5 = @("Block Bootstrapping", 0, 0, 8, 0, 3.108082487074975)
    6 = @("LSPMW", 0, 2, 6, 0, 1.611972106168202)
    7 = @("AREPD", 0, 0, 8, 0, 2.937753049473157)
    8 = @("MCPS", 0, 0, 8, 0, 1.187908763415526)
    9 = @("DeepAR", 0, 0, 8, 0, 1.126506566211152)
    10 = @("EnCQR-LSTM", 0, 3, 5, 0, 2.297252780636199)
}
foreach ($r in $resumenRows.Keys) {
    $vals = $resumenRows[$r]
    $wsResumen.Range("A$r").Value = $vals[0]
    $wsResumen.Range("B$r").Value = $vals[1]
    $wsResumen.Range("C$r").Value = $vals[2]
    $wsResumen.Range("D$r").Value = $vals[3]
    $wsResumen.Range("E$r").Value = $vals[4]
    $wsResumen.Range("F$r").Value = $vals[5]
}

$wb.Save()
